# Auto-generated edit script applying scheduled-runner price/profit updates
# to the Seraph_Profits sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 116.4
$ws.Range("I11").Value = 116.4
$ws.Range("K11").Value = 116.4
$ws.Range("M11").Value = 23.59999999999999
$ws.Range("H40").Value = 2545.2727
$ws.Range("I40").Value = 2166.3333
$ws.Range("J40").Value = 3000
$ws.Range("K40").Value = 2166.3333
$ws.Range("L40").Value = 3000
$ws.Range("M40").Value = -1991.3333
$ws.Range("N40").Value = -3350
$ws.Range("H42").Value = 73.818184
$ws.Range("I42").Value = 22.833334
$ws.Range("K42").Value = 68.50000199999999
$ws.Range("M42").Value = 161.499998
$ws.Range("H55").Value = 1526.4117
$ws.Range("J55").Value = 2885.5715
$ws.Range("L55").Value = 2885.5715
$ws.Range("N55").Value = -3313.5715
$ws.Range("H88").Value = 3763.3635
$ws.Range("J88").Value = 3711.111
$ws.Range("L88").Value = 3711.111
$ws.Range("N88").Value = -4523.111
$ws.Range("H91").Value = 3763.3635
$ws.Range("J91").Value = 3711.111
$ws.Range("L91").Value = 3711.111
$ws.Range("N91").Value = -6519.111
$ws.Range("H106").Value = 30616.883
$ws.Range("I106").Value = 32807.46
$ws.Range("J106").Value = 23497.5
$ws.Range("K106").Value = 32807.46
$ws.Range("L106").Value = 23497.5
$ws.Range("M106").Value = -32176.46
$ws.Range("N106").Value = -24759.5
$ws.Range("H111").Value = 1077.1
$ws.Range("I111").Value = 1077.1
$ws.Range("J111").Value = 0
$ws.Range("K111").Value = 3231.3
$ws.Range("L111").Value = 0
$ws.Range("M111").ClearContents()
$ws.Range("N111").Value = -164.2999999999997
$ws.Range("H129").Value = 1813.5454
$ws.Range("I129").Value = 564.2857
$ws.Range("K129").Value = 1692.8571
$ws.Range("M129").Value = 3307.1429

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 157.875
$ws.Range("I5").Value = 166.14285
$ws.Range("K5").Value = 166.14285
$ws.Range("M5").Value = -54.14285000000001
$ws.Range("H32").Value = 7730.316
$ws.Range("I32").Value = 5436.1562
$ws.Range("J32").Value = 19965.834
$ws.Range("K32").Value = 5436.1562
$ws.Range("L32").Value = 19965.834
$ws.Range("M32").Value = -5149.1562
$ws.Range("N32").Value = -20539.834
$ws.Range("H45").Value = 4666.3335
$ws.Range("I45").Value = 4499.5
$ws.Range("K45").Value = 4499.5
$ws.Range("M45").Value = -4122.5
$ws.Range("H53").Value = 0
$ws.Range("J53").Value = 0
$ws.Range("L53").ClearContents()
$ws.Range("N53").Value = 0
$ws.Range("H88").Value = 1535
$ws.Range("I88").Value = 797
$ws.Range("J88").Value = 1977.8
$ws.Range("K88").Value = 797
$ws.Range("L88").Value = 1977.8
$ws.Range("M88").Value = -391
$ws.Range("N88").Value = -2789.8
$ws.Range("H91").Value = 1535
$ws.Range("I91").Value = 797
$ws.Range("J91").Value = 1977.8
$ws.Range("K91").Value = 797
$ws.Range("L91").Value = 1977.8
$ws.Range("M91").Value = 607
$ws.Range("N91").Value = -4785.8
$ws.Range("H102").Value = 2434.75
$ws.Range("I102").Value = 870
$ws.Range("K102").Value = 870
$ws.Range("M102").Value = 752
$ws.Range("H122").Value = 669219.9399999999
$ws.Range("I122").Value = 716628.5
$ws.Range("K122").Value = 2149885.5
$ws.Range("M122").Value = -2147435.5

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 157.875
$ws.Range("I4").Value = 166.14285
$ws.Range("K4").Value = 166.14285
$ws.Range("M4").Value = -51.14285000000001
$ws.Range("H32").Value = 3000
$ws.Range("J32").Value = 3000
$ws.Range("L32").Value = 3000
$ws.Range("N32").Value = -3768
$ws.Range("H37").Value = 350
$ws.Range("I37").Value = 350
$ws.Range("K37").Value = 350
$ws.Range("M37").Value = -213
$ws.Range("H94").Value = 1352.5333
$ws.Range("I94").Value = 1099.0769
$ws.Range("J94").Value = 3000
$ws.Range("K94").Value = 1099.0769
$ws.Range("L94").Value = 3000
$ws.Range("M94").Value = -648.0769
$ws.Range("N94").Value = -3902
$ws.Range("H105").Value = 3687.6667
$ws.Range("I105").Value = 3586.125
$ws.Range("J105").Value = 4500
$ws.Range("K105").Value = 3586.125
$ws.Range("L105").Value = 4500
$ws.Range("M105").Value = -1839.125
$ws.Range("N105").Value = -7994
$ws.Range("H138").Value = 99999.5
$ws.Range("J138").Value = 99999.5
$ws.Range("L138").Value = 99999.5
$ws.Range("N138").Value = -110279.5
$ws.Range("H141").Value = 75000
$ws.Range("J141").Value = 75000
$ws.Range("L141").Value = 75000
$ws.Range("N141").Value = -85360

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 5099.9165
$ws.Range("J4").Value = 5649.875
$ws.Range("L4").Value = 5649.875
$ws.Range("N4").Value = -5873.875
$ws.Range("H16").Value = 3749
$ws.Range("I16").Value = 3749
$ws.Range("K16").Value = 3749
$ws.Range("M16").Value = -3462
$ws.Range("H31").Value = 3021.6553
$ws.Range("I31").Value = 2640.72
$ws.Range("J31").Value = 5402.5
$ws.Range("K31").Value = 2640.72
$ws.Range("L31").Value = 5402.5
$ws.Range("M31").Value = -2345.72
$ws.Range("N31").Value = -5992.5
$ws.Range("H34").Value = 3021.6553
$ws.Range("I34").Value = 2640.72
$ws.Range("J34").Value = 5402.5
$ws.Range("K34").Value = 2640.72
$ws.Range("L34").Value = 5402.5
$ws.Range("M34").Value = -2438.72
$ws.Range("N34").Value = -5806.5
$ws.Range("H99").Value = 14478.523
$ws.Range("I99").Value = 12681.286
$ws.Range("J99").Value = 15377.143
$ws.Range("K99").Value = 12681.286
$ws.Range("L99").Value = 15377.143
$ws.Range("M99").Value = -11183.286
$ws.Range("N99").Value = -18373.143
$ws.Range("H113").Value = 3749
$ws.Range("I113").Value = 3749
$ws.Range("K113").Value = 3749
$ws.Range("M113").Value = -1579
$ws.Range("H126").Value = 14478.523
$ws.Range("I126").Value = 12681.286
$ws.Range("J126").Value = 15377.143
$ws.Range("K126").Value = 38043.858
$ws.Range("L126").Value = 46131.429
$ws.Range("M126").Value = -35573.858
$ws.Range("N126").Value = -51071.429
$ws.Range("H134").Value = 1713.8158
$ws.Range("I134").Value = 1513.4375
$ws.Range("J134").Value = 2782.5
$ws.Range("K134").Value = 4540.3125
$ws.Range("L134").Value = 8347.5
$ws.Range("M134").Value = -2005.3125
$ws.Range("N134").Value = -13417.5

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1108.1538
$ws.Range("J131").Value = 1241.25
$ws.Range("L131").Value = 3723.75
$ws.Range("N131").Value = -13803.75

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H19").Value = 122.85714
$ws.Range("I19").Value = 122.85714
$ws.Range("K19").Value = 122.85714
$ws.Range("M19").Value = 165.14286
$ws.Range("H97").Value = 0
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 0
$ws.Range("L97").ClearContents()
$ws.Range("M97").ClearContents()
$ws.Range("N97").Value = 0
$ws.Range("H99").Value = 9077.429
$ws.Range("J99").Value = 27250
$ws.Range("L99").Value = 27250
$ws.Range("N99").Value = -31742
$ws.Range("H102").Value = 2740.3333
$ws.Range("I102").Value = 2740.3333
$ws.Range("K102").Value = 2740.3333
$ws.Range("M102").Value = -1118.3333

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4468.8
$ws.Range("I7").Value = 4468.8
$ws.Range("K7").Value = 4468.8
$ws.Range("M7").Value = -4356.8
$ws.Range("H46").Value = 5241.5
$ws.Range("I46").Value = 2248.75
$ws.Range("J46").Value = 7236.6665
$ws.Range("K46").Value = 2248.75
$ws.Range("L46").Value = 7236.6665
$ws.Range("M46").Value = -2060.75
$ws.Range("N46").Value = -7612.6665
$ws.Range("H61").Value = 3784
$ws.Range("I61").Value = 3752.4285
$ws.Range("J61").Value = 4005
$ws.Range("K61").Value = 3752.4285
$ws.Range("L61").Value = 4005
$ws.Range("M61").Value = -3550.4285
$ws.Range("N61").Value = -4409
$ws.Range("H87").Value = 17500
$ws.Range("J87").Value = 17500
$ws.Range("L87").Value = 17500
$ws.Range("N87").Value = -19746
$ws.Range("H90").Value = 17500
$ws.Range("J90").Value = 17500
$ws.Range("L90").Value = 52500
$ws.Range("N90").Value = -63732
$ws.Range("H93").Value = 1861.375
$ws.Range("I93").Value = 1770.1428
$ws.Range("K93").Value = 1770.1428
$ws.Range("M93").Value = -522.1428000000001
$ws.Range("H113").Value = 3784
$ws.Range("I113").Value = 3752.4285
$ws.Range("J113").Value = 4005
$ws.Range("K113").Value = 3752.4285
$ws.Range("L113").Value = 4005
$ws.Range("M113").Value = -1582.4285
$ws.Range("N113").Value = -8345
$ws.Range("H126").Value = 4468.8
$ws.Range("I126").Value = 4468.8
$ws.Range("K126").Value = 13406.4
$ws.Range("M126").Value = -10936.4

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2998.5715
$ws.Range("I122").Value = 2999
$ws.Range("J122").Value = 2997.5
$ws.Range("K122").Value = 8997
$ws.Range("L122").Value = 8992.5
$ws.Range("M122").Value = -6547
$ws.Range("N122").Value = -13892.5

